$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for the new points table
$ws.Range("N3").Value = "p1"
$ws.Range("O3").Value = "p2"
$ws.Range("P3").Value = "p3"

# Row 4 values + average
$ws.Range("N4").Value = -2.5
$ws.Range("O4").Value = -2.5
$ws.Range("P4").Value = 2
$ws.Range("Q4").Formula = "=AVERAGE(N4:P4)"

# Row 5 values
$ws.Range("P5").Value = 2

# Row 6 values
$ws.Range("N6").Value = -2.5
$ws.Range("O6").Value = 2.5
$ws.Range("P6").Value = 2

# Row 5/6 average as one shared formula (relative refs shift per row)
$ws.Range("Q5:Q6").Formula = "=AVERAGE(N5:P5)"

# Update selection to mirror final state from the diff
$ws.Range("O5").Select()
